$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-14 10:23:32"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
